$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 167884
$ws.Range("C4").Value = 158764
$ws.Range("C7").Value = 5.43
$ws.Range("C8").Value = 65.33
